$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue 'D2' '292.61'
Set-TextValue 'E2' '-0.35%'
Set-TextValue 'D3' '40.41'
Set-TextValue 'E3' '1.08%'
Set-TextValue 'E4' '-0.44%'
Set-TextValue 'B6' 'FTXToken'
Set-TextValue 'C6' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D6' '1.569'
Set-TextValue 'E6' '1.49%'
Set-TextValue 'B7' 'MXToken'
Set-TextValue 'C7' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D7' '0.9237'
Set-TextValue 'E7' '0.00%'
Set-TextValue 'B8' 'BTSEToken'
Set-TextValue 'C8' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue 'D8' '2.381'
Set-TextValue 'E8' '-0.76%'
Set-TextValue 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D9' '0.1185'
Set-TextValue 'E9' '-0.20%'
Set-TextValue 'B10' 'WazirX'
Set-TextValue 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1810'
Set-TextValue 'E10' '3.12%'
Set-TextValue 'B11' 'BitrueCoin'
Set-TextValue 'C11' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D11' '0.04387'
Set-TextValue 'E11' '5.31%'
Set-TextValue 'B12' 'MandalaExchangeToken'
Set-TextValue 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.08769'
Set-TextValue 'E12' '1.25%'
Set-TextValue 'B13' 'BitMartToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D13' '0.1053'
Set-TextValue 'E13' '-0.12%'
Set-TextValue 'B14' 'TigerCash'
Set-TextValue 'C14' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D14' '0.005787'
Set-TextValue 'E14' '0.27%'
Set-TextValue 'B15' 'LEO'
Set-TextValue 'C15' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D15' '3.340'
Set-TextValue 'E15' '-1.10%'
Set-TextValue 'B16' 'GateToken'
Set-TextValue 'C16' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D16' '4.288'
Set-TextValue 'E16' '-0.40%'
Set-TextValue 'B17' 'BitpandaEcosystemToken'
Set-TextValue 'C17' 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue 'D17' '0.3318'
Set-TextValue 'E17' '0.70%'
Set-TextValue 'B18' 'MCDex'
Set-TextValue 'C18' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D18' '7.846'
Set-TextValue 'E18' '3.43%'
Set-TextValue 'B19' 'ProBitToken'
Set-TextValue 'C19' 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue 'D19' '0.1391'
Set-TextValue 'E19' '3.56%'
Set-TextValue 'B20' 'ZBToken'
Set-TextValue 'C20' 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue 'D20' '0.2957'
Set-TextValue 'E20' '5.38%'
Set-TextValue 'B21' 'BitForexToken'
Set-TextValue 'C21' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D21' '0.001267'
Set-TextValue 'E21' '-0.35%'
Set-TextValue 'E22' '2.84%'
Set-TextValue 'D23' '0.001259'
Set-TextValue 'E23' '-1.89%'
Set-TextValue 'D24' '0.003737'
Set-TextValue 'E24' '-4.28%'
Set-TextValue 'D25' '0.0001250'
Set-TextValue 'E25' '-3.47%'
Set-TextValue 'D26' '0.0003722'
Set-TextValue 'E26' '-0.22%'
Set-TextValue 'D38' '0.02347'
Set-TextValue 'E38' '1.56%'
Set-TextValue 'D39' '0.05072'
Set-TextValue 'E39' '1.42%'
Set-TextValue 'D40' '0.005790'
Set-TextValue 'E40' '34.36%'
Set-TextValue 'D41' '0.007800'
Set-TextValue 'E41' '0.94%'
Set-TextValue 'D42' '0.1290'
Set-TextValue 'E42' '1.11%'
Set-TextValue 'D43' '0.007384'
Set-TextValue 'D44' '0.008039'
Set-TextValue 'E44' '15.32%'
Set-TextValue 'D45' '0.2914'
Set-TextValue 'E45' '-8.63%'
Set-TextValue 'E46' '-3.95%'
Set-TextValue 'E47' '-0.21%'
Set-TextValue 'D49' '0.004200'
Set-TextValue 'E49' '-0.22%'
Set-TextValue 'D50' '0.00002100'
Set-TextValue 'E50' '-0.21%'
Set-TextValue 'D51' '0.0002000'
Set-TextValue 'E51' '-0.21%'
